# Updates the cryptos price/volume snapshot (GitHub Actions scraper refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.664.71"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "3.432.14"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.71"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.42"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "3.432.39"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.130"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.90"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.407"
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("D13").Value = "4.028.34"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.48"
$ws.Range("E15").Value = "  -7.73%  "
$ws.Range("D16").Value = "65.750.25"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000169"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "3.435.46"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.90"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.71"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.35"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.52"
$ws.Range("E22").Value = "  -2.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.99"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.528"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000120"
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.62"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.176"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.51"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.68"
$ws.Range("E32").Value = "  -3.30%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  -5.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.94"
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.17"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.877"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.43"
$ws.Range("E39").Value = "  +5.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.59"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.75"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("D42").Value = "2.755.97"
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.43"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.41"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0675"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.15"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.41"
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0287"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "324.59"
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.21"
$ws.Range("E51").Value = "  +0.40%  "
